# Junction_Flooding_425.xlsx refresh: new sampling window ("custom accuracy + 1000 data points")
# - replace the 4 data rows (rows 2-5) with a new time-window of readings
# - drop the old 5th data row (sheet row 6), shrinking the used range to A1:AH5
# - a handful of columns grow by one character of width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New readings for rows 2-5 (columns A:AH == Time, J1..J33)
# ---------------------------------------------------------------------------
$data = @(
    @(45106.50694444445, 4.96, 4.926, 1.074, 9.69, 9.442, 3.413, 8.191000000000001, 4.526, 2.102, 4.543, 5.171, 3.55, 0.877, 3.454, 4.766, 2.485, 1.205, 0, 48.94, 9.632, 2.995, 5.545, 4.641, 0.661, 7.552, 1.905, 4.889, 2.236, 4.485, 0.29, 5.812, 2.376, 3.752),
    @(45106.51388888889, 13.183, 10.227, 0.82, 27.862, 23.8, 9.964, 34.452, 15.134, 6.943, 10.879, 11.591, 11.731, 3.163, 9.955, 14.332, 8.27, 0.698, 0.194, 146.788, 27.924, 9.087999999999999, 18.702, 10.478, 1.354, 18.899, 7.655, 8.295, 8.156000000000001, 12.04, 0.033, 30.582, 5.609, 11.403),
    @(45106.52083333334, 6.894, 5.339, 0.469, 14.4, 12.427, 5.114, 22.625, 7.788, 3.639, 5.691, 6, 6.09, 1.649, 5.118, 7.482, 4.301, 0.459, 0.045, 72.566, 14.577, 4.659, 9.805, 5.499, 0.699, 11.722, 3.869, 4.444, 4.185, 6.283, 0, 20.425, 2.962, 5.863),
    @(45106.52777777778, 20.31, 15.36, 0.85, 43.72, 36.54, 15.74, 58.34, 24.18, 10.95, 16.56, 17.69, 18.56, 5.04, 15.67, 22.53, 13.05, 0.48, 0.47, 232.84, 43.85, 14.42, 29.68, 15.87, 2.07, 29.5, 12.55, 11.8, 13.11, 18.61, 0, 52.53, 8.460000000000001, 18.06)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the now-obsolete last data row (sheet row 6); this also shrinks
#    the sheet dimension from A1:AH6 down to A1:AH5 automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3) Widen a subset of columns by one character.
#    Note: Excel's COM `ColumnWidth` is character-width based and is offset
#    from the raw OOXML `width` (stored width) by the standard 5px padding
#    (5/6 of a character at the default Calibri 11 / MDW=7 font metrics).
#    So to land on an integer stored width of N we set ColumnWidth = N - 5/6.
# ---------------------------------------------------------------------------
$padding = 5 / 6

$columnsToWidthEight = @(2, 3, 6, 9, 11, 12, 13, 16, 23, 24, 26, 34)
foreach ($col in $columnsToWidthEight) {
    $ws.Columns.Item($col).ColumnWidth = 8 - $padding
}

$ws.Columns.Item(20).ColumnWidth = 9 - $padding
